$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- The credential table had its 4th row ("test121@test.com" / "test123")
# removed, so every row below it shifted up by one, and a brand-new row of
# credentials ("testselenium12345@gmail.com" / "test123") was appended at
# what is now row 7. Re-create that by shifting the surviving rows up and
# writing the new pair into row 7.

$a5 = $ws.Range("A5").Value()
$b5 = $ws.Range("B5").Value()
$a6 = $ws.Range("A6").Value()
$b6 = $ws.Range("B6").Value()
$a7 = $ws.Range("A7").Value()
$b7 = $ws.Range("B7").Value()

$ws.Range("A4").Value = $a5
$ws.Range("B4").Value = $b5
$ws.Range("A5").Value = $a6
$ws.Range("B5").Value = $b6
$ws.Range("A6").Value = $a7
$ws.Range("B6").Value = $b7

$ws.Range("A7").Value = "testselenium12345@gmail.com"
$ws.Range("B7").Value = "test123"

# New row 7 height matches what the sheet re-computed for rows 4-7 after
# the edit (same as the existing hyperlink row 2).
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 13.8
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8

# A7 now holds an email address, so it becomes a mailto hyperlink just like
# A2 already is.
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:testselenium12345@gmail.com", "", "", "testselenium12345@gmail.com")

# Give A7/B7 the same look as the rest of the table: A7 picks up the
# hyperlink-cell formatting already used by A2, B7 goes back to the sheet's
# plain/default formatting (copied from an untouched, default-styled cell).
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Hyperlinks.Add registers its own builtin "Hyperlink" cell style even
# though we immediately override the cell's actual formatting above; drop
# the now-unused style definition again.
$wb.Styles.Item("Hyperlink").Delete()

# Final cursor position left by the edit.
$null = $ws.Range("B13").Select()
